$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$deptWs = $wb.Worksheets.Item(2)
$titleWs = $wb.Worksheets.Item(3)

# Row 2 - fix: employee code left blank, names/email corrected
$ws.Range("A2").ClearContents()
$ws.Range("B2").Value = "Ahmedov1"
$ws.Range("C2").Value = "Suhrob1"
$ws.Range("E2").Value = "EMPLOYEE"
$ws.Range("G2").Value = "2-stage"
$ws.Range("H2").Value = "CEO"
$ws.Range("J2").Value = "suhrob1@gmail.com"

# Row 3
$ws.Range("A3").Value = "DK0002"
$ws.Range("B3").Value = "Ahmedov2"
$ws.Range("C3").Value = "Suhrob2"
$ws.Range("E3").Value = "EMPLOYEE"
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = "2-stage"
$ws.Range("H3").Value = "CEO"
$ws.Range("I3").Value = "Kjhjgf"
$ws.Range("J3").Value = "suhrob2@gmail.com"

# Row 4
$ws.Range("A4").Value = "DK0003"
$ws.Range("B4").Value = "Ahmedov3"
$ws.Range("C4").Value = "Suhrob3"
$ws.Range("E4").Value = "EMPLOYEE"
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = "2-stage"
$ws.Range("H4").Value = "CEO"
$ws.Range("J4").Value = "suhrob3gmail.com"

# Row 5
$ws.Range("A5").Value = "DK0004"
$ws.Range("B5").Value = "Ahmedov4"
$ws.Range("C5").Value = "Suhrob4"
$ws.Range("E5").Value = "EMPLOYEE"
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = "2-stage"
$ws.Range("H5").Value = "CEO"
$ws.Range("J5").Value = "suhrob4@gmail.com"

# Row 6
$ws.Range("A6").Value = "DK0005"
$ws.Range("B6").Value = "Ahmedov5"
$ws.Range("C6").Value = "Suhrob5"
$ws.Range("E6").Value = "EMPLOYEE"
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = "2-stage"
$ws.Range("H6").Value = "CEO"
$ws.Range("J6").Value = "suhrob4@gmail.com"

# Row 7
$ws.Range("A7").Value = "１１１ｋ"
$ws.Range("B7").Value = "Ahmedov6"
$ws.Range("C7").Value = "Suhrob6"
$ws.Range("E7").Value = "EMPLOYEE"
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = "2-stage"
$ws.Range("H7").Value = "CEO"
$ws.Range("J7").Value = "suhrob6@gmail.com"

# Row 8
$ws.Range("A8").Value = "DK0007"
$ws.Range("B8").Value = "  "
$ws.Range("C8").Value = "Suhrob7"
$ws.Range("E8").Value = "EMPLOYEE"
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = "2-stage"
$ws.Range("H8").Value = "CEO"
$ws.Range("J8").Value = "suhrob7@gmail.com"

# Row 9
$ws.Range("A9").Value = "   "
$ws.Range("B9").Value = "Ahmedov8"
$ws.Range("C9").Value = "Suhrob8"
$ws.Range("E9").Value = "EMPLOYEE"
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = "2-stage"
$ws.Range("H9").Value = "CEO"
$ws.Range("J9").Value = "suhrob8@gmail.com"

# Row 10
$ws.Range("A10").Value = "DK0009"
$ws.Range("B10").Value = "Ahmedov9"
$ws.Range("C10").Value = "Suhrob9"
$ws.Range("E10").Value = "EMPLOYEE"
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = "2-stage"
$ws.Range("H10").Value = "CEO"
$ws.Range("J10").Value = "suhrob9@gmail.com"

# Dropdown source lists (部署/役職 sheets) - re-entered values
$deptWs.Range("A2").Value = "1-stage"
$titleWs.Range("A2").Value = "Software Engineer"

$ws.Activate()
$ws.Range("B14").Select()
